# Machine Learning Bootcamp Sciprog.pptx - apply commit "pushing dropbox to git"
#
# Summary of the edit (reconstructed from the sldId-based slide identity,
# which is stable across the OOXML part renames that a naive file-diff
# shows): three brand-new slides were inserted into the deck, while every
# pre-existing slide kept its original content and relative order:
#
#   1. Title slide                              (unchanged)
#   2. Topics                                    (unchanged)
#   3. [NEW] Multivariate regression
#   4. What is Machine Learning / Algorithms...  (unchanged)
#   5. What is Machine Learning / Terminology    (unchanged)
#   6. Types of Machine Learning                 (unchanged)
#   7. Regression (+ picture)                    (unchanged)
#   8. [NEW] Regression / scikit-learn
#   9. Regression: Overfitting (+ picture)        (unchanged)
#   10. [NEW] Mission
#   11. Classification (+ picture)                (unchanged)
#   12. (blank placeholder slide)                 (unchanged)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) New slide "Multivariate regression", inserted right after the
#    "Topics" slide (position 2 -> new slide becomes position 3).
# ---------------------------------------------------------------------
$s1 = $p.Slides.Add(3, 2)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Multivariate regression"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Both of us look at this"

# ---------------------------------------------------------------------
# 2) New slide "Regression" / scikit-learn how-to, inserted right after
#    the "Regression" (+ picture) slide, which is now at position 7.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(8, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Regression"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "How to do in sci kit learn `rLinear and polynomial regression"

# ---------------------------------------------------------------------
# 3) New slide "Mission", inserted right after the "Regression:
#    Overfitting" (+ picture) slide, which is now at position 9.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(10, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Mission"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "In this workshop we will give you a quick demo into Machine Learning using Python"
